$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11; existing rows 11-39 shift down to 12-40.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly price entry.
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C11").Value = "Ñuble"
$ws.Range("D11").Value = "2022-07-22"
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 100112013
$ws.Range("G11").Value = "Alcachofa"
$ws.Range("H11").Value = "Argentina(o)"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 60
$ws.Range("K11").Value = 15000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 15500
$ws.Range("N11").Value = "`$/caja 50 unidades"
$ws.Range("O11").Value = "Provincia de Limarí"
$ws.Range("P11").Value = 310
$ws.Range("Q11").Value = 50
$ws.Range("R11").Value = "Hortaliza"
